$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "O3" = 1.25
    "P3" = 3.75
    "Q3" = 1.9
    "R3" = 1.95
    "AB6" = 26
    "AC6" = 10
    "AD6" = 13.5
    "AE6" = 25
    "AF6" = 100
    "AG6" = 37
    "AI6" = 37
    "AL6" = 110
    "AM6" = 700
    "AN6" = 3.2
    "AO6" = 5
    "AP6" = 14.5
    "AQ6" = 11
    "AR6" = 32
    "AS6" = 175
    "AT6" = 3.6
    "AU6" = 9.25
    "AV6" = 80
    "AW6" = 11.75
    "AX6" = 70
    "AY6" = 55
    "BA6" = 450
    "G6" = 1.19
    "H6" = 5.9
    "I6" = 11.25
    "K6" = 2.8
    "L6" = 9
    "O6" = 1.12
    "P6" = 5.4
    "Q6" = 1.38
    "R6" = 2.82
    "T6" = 3.6
    "U6" = 1.91
    "V6" = 1.8
    "W6" = 10
    "Y6" = 9.75
    "Z6" = 7.6
    "AA7" = 50
    "AC7" = 7.2
    "AD7" = 7
    "AE7" = 16.5
    "AG7" = 6.3
    "AH7" = 7.3
    "AI7" = 8
    "AJ7" = 12
    "AK7" = 13.5
    "AN7" = 6.7
    "AO7" = 29
    "AQ7" = 175
    "AR7" = 200
    "AS7" = 450
    "AT7" = 2.7
    "AU7" = 7.6
    "AW7" = 3.45
    "AX7" = 8.25
    "AY7" = 18.5
    "AZ7" = 27
    "BA7" = 65
    "G7" = 5
    "H7" = 3.5
    "I7" = 1.62
    "J7" = 5.2
    "K7" = 2.15
    "L7" = 2.22
    "N7" = 7.2
    "O7" = 1.31
    "P7" = 3.15
    "Q7" = 1.93
    "R7" = 1.8
    "T7" = 2.7
    "V7" = 1.82
    "Y7" = 16
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "done"